$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("5").Delete()
$ws.Rows("3").Delete()

$ws.Range("A3").Select()
